# Update RandomForestClassifier_tpesearch-val_df results for
# target_col == previous_concussions (first finished training run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("0")
$ws.Range("B7").Value = 0.3333333333333333
$ws.Range("C7").Value = 0.4
$ws.Range("D7").Value = 0.3636363636363636

# Row 8 ("1")
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 0.4285714285714285
$ws.Range("D8").Value = 0.4615384615384615

# Row 9 ("accuracy")
$ws.Range("B9").Value = 0.4166666666666667
$ws.Range("C9").Value = 0.4166666666666667
$ws.Range("D9").Value = 0.4166666666666667
$ws.Range("E9").Value = 0.4166666666666667

# Row 10 ("macro avg")
$ws.Range("B10").Value = 0.4166666666666666
$ws.Range("C10").Value = 0.4142857142857143
$ws.Range("D10").Value = 0.4125874125874125

# Row 11 ("weighted avg")
$ws.Range("B11").Value = 0.4305555555555555
$ws.Range("C11").Value = 0.4166666666666667
$ws.Range("D11").Value = 0.4207459207459207

# Row 17 ("0")
$ws.Range("B17").Value = 0.4545454545454545
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = 0.4761904761904762

# Row 18 ("1")
$ws.Range("B18").Value = 0.6153846153846154
$ws.Range("D18").Value = 0.5925925925925927

# Row 19 ("accuracy")
$ws.Range("B19").Value = 0.5416666666666666
$ws.Range("C19").Value = 0.5416666666666666
$ws.Range("D19").Value = 0.5416666666666666
$ws.Range("E19").Value = 0.5416666666666666

# Row 20 ("macro avg")
$ws.Range("B20").Value = 0.534965034965035
$ws.Range("C20").Value = 0.5357142857142857
$ws.Range("D20").Value = 0.5343915343915344

# Row 21 ("weighted avg")
$ws.Range("B21").Value = 0.5483682983682984
$ws.Range("C21").Value = 0.5416666666666666
$ws.Range("D21").Value = 0.5440917107583775
